# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and bump the
# related handoff/handback timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-23 00:34:55"

# zh-cn sheet: ... | Status(C) | Latest Handoff File(D) | Latest Handoff Datetime(E) | ...
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-23 00:34:51"

# de-de sheet: ... | Status(C) | Latest Handoff File(D) | Latest Handoff Datetime(E) | ...
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-23 00:34:55"
